$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two old 2018 events that occupied rows 7 and 8:
#   row 7: "Mechanisms involved in the prevention of type 1 diabetes onset by Lactobacillus johnsonii N6.2" (Webinar)
#   row 8: "NCI Containers and Workflows Interest Group Seminar" (Seminar)
# Row 8's url cell (G8) carries a hyperlink; drop it explicitly first so it
# doesn't linger as an orphaned relationship once the row is gone.
foreach ($h in $ws.Hyperlinks) {
  if ($h.Range.Address() -eq '$G$8') {
    $h.Delete()
  }
}

# Deleting these whole rows shifts the remaining "Bioinformatics User Forum Meeting" row up to become row 7.
$ws.Range("A7:A8").EntireRow.Delete()

# Restore the selection to what it was after the manual row-delete
# (whole rows 7:8 selected, landing on A7 once they're removed).
$null = $ws.Range("A7:XFD8").Select()
